$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.25
$ws.Range("D2").Value = "7:45pm"

$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A3").Value = 43828
$ws.Range("B3").Value = 2.5
$ws.Range("C3").Value = "7:44pm"
$ws.Range("D3").Value = "10:25pm"

$ws.Range("F7").Select()
